# Attendance sheet update: Google Sheets / email-notification export changes.
# Source sheet originally tracked 2025-11-07 .. 2025-11-21 with columns:
#   A=Roll No, B=Name, C..Q=dates, R=Present, S=Total, T=Attendance %
# Target sheet adds an Email column, two placeholder "2025-11-22_x/_y" columns,
# a genuine 2025-11-22 date column, a newly-enrolled student (Abhishek Pathak),
# and fills in the previously-blank attendance row for Shubham Pitekar.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural edits -----------------------------------------------------

# Make room for two new columns right after "Roll No"/"Name"/old-date column C
# (they become the "2025-11-22_x" / "2025-11-22_y" columns; column C itself is
# repurposed from a date column into the new "Email" column).
$ws.Columns("D:E").Insert()

# Make room for the real "2025-11-22" date column just before the old
# "Present" column (which has since shifted from R to T).
$ws.Columns("T:T").Insert()

# Make room for the newly-enrolled student "Abhishek Pathak" (EC4226), who is
# inserted right after "Anushka Mote" (row 4) and before "Vaishnavi Pawar".
$ws.Rows("5:5").Insert()

# --- Cell content -----------------------------------------------------------

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W")
$numericCols = @("U","V","W")

$data = @(
    @{ "A" = "Roll No"; "B" = "Name"; "C" = "Email"; "D" = "2025-11-22_x"; "E" = "2025-11-22_y"; "F" = "2025-11-08"; "G" = "2025-11-09"; "H" = "2025-11-10"; "I" = "2025-11-11"; "J" = "2025-11-12"; "K" = "2025-11-13"; "L" = "2025-11-14"; "M" = "2025-11-15"; "N" = "2025-11-16"; "O" = "2025-11-17"; "P" = "2025-11-18"; "Q" = "2025-11-19"; "R" = "2025-11-20"; "S" = "2025-11-21"; "T" = "2025-11-22"; "U" = "Present"; "V" = "Total"; "W" = "Attendance %" };
    @{ "A" = "EC4202"; "B" = "Divya Bhagas"; "C" = "divyabhagas@gmail.com"; "D" = "❌"; "E" = "❌"; "F" = "❌"; "G" = "❌"; "H" = "❌"; "I" = "❌"; "J" = "❌"; "K" = "❌"; "L" = "❌"; "M" = "❌"; "N" = "❌"; "O" = "❌"; "P" = "❌"; "Q" = "❌"; "R" = "❌"; "S" = "❌"; "T" = "❌"; "U" = 0; "V" = 17; "W" = 0 };
    @{ "A" = "EC4206"; "B" = "Sai Kenekar"; "C" = "saikanekar@gmail.com"; "D" = "❌"; "E" = "❌"; "F" = "❌"; "G" = "❌"; "H" = "❌"; "I" = "❌"; "J" = "❌"; "K" = "❌"; "L" = "❌"; "M" = "❌"; "N" = "❌"; "O" = "❌"; "P" = "❌"; "Q" = "❌"; "R" = "❌"; "S" = "❌"; "T" = "❌"; "U" = 0; "V" = 17; "W" = 0 };
    @{ "A" = "EC4220"; "B" = "Anushka Mote"; "C" = "anushska484m@gmail.com"; "D" = "❌"; "E" = "❌"; "F" = "❌"; "G" = "❌"; "H" = "❌"; "I" = "❌"; "J" = "❌"; "K" = "❌"; "L" = "❌"; "M" = "❌"; "N" = "❌"; "O" = "❌"; "P" = "❌"; "Q" = "✅"; "R" = "❌"; "S" = "❌"; "T" = "❌"; "U" = 1; "V" = 17; "W" = 5.9 };
    @{ "A" = "EC4226"; "B" = "Abhishek Pathak"; "C" = "abhipathak2513@gmail.com"; "D" = "❌"; "E" = "❌"; "F" = "❌"; "G" = "❌"; "H" = "❌"; "I" = "❌"; "J" = "❌"; "K" = "❌"; "L" = "❌"; "M" = "❌"; "N" = "❌"; "O" = "❌"; "P" = "❌"; "Q" = "❌"; "R" = "❌"; "S" = "❌"; "T" = "❌"; "U" = 0; "V" = 17; "W" = 0 };
    @{ "A" = "EC4231"; "B" = "Vaishnavi Pawar"; "C" = "vaishnavipawar@gmail.com"; "D" = "❌"; "E" = "❌"; "F" = "❌"; "G" = "❌"; "H" = "❌"; "I" = "❌"; "J" = "❌"; "K" = "❌"; "L" = "❌"; "M" = "❌"; "N" = "❌"; "O" = "❌"; "P" = "❌"; "Q" = "❌"; "R" = "❌"; "S" = "❌"; "T" = "❌"; "U" = 0; "V" = 17; "W" = 0 };
    @{ "A" = "EC4233"; "B" = "Sagar Pawar"; "C" = "sagarpawar@gmail.com"; "D" = "❌"; "E" = "❌"; "F" = "❌"; "G" = "❌"; "H" = "❌"; "I" = "❌"; "J" = "❌"; "K" = "❌"; "L" = "❌"; "M" = "❌"; "N" = "❌"; "O" = "❌"; "P" = "❌"; "Q" = "❌"; "R" = "❌"; "S" = "❌"; "T" = "❌"; "U" = 0; "V" = 17; "W" = 0 };
    @{ "A" = "EC4236"; "B" = "Shubham Phad"; "C" = "shubhamphad03@gmail.com"; "D" = "❌"; "E" = "❌"; "F" = "❌"; "G" = "❌"; "H" = "❌"; "I" = "❌"; "J" = "❌"; "K" = "❌"; "L" = "❌"; "M" = "❌"; "N" = "❌"; "O" = "❌"; "P" = "❌"; "Q" = "❌"; "R" = "❌"; "S" = "❌"; "T" = "❌"; "U" = 0; "V" = 17; "W" = 0 };
    @{ "A" = "EC4237"; "B" = "Shubham Pitekar"; "C" = "shubhampitekar2323@gmail.com"; "D" = "❌"; "E" = "❌"; "F" = "❌"; "G" = "❌"; "H" = "❌"; "I" = "❌"; "J" = "❌"; "K" = "❌"; "L" = "❌"; "M" = "✅"; "N" = "❌"; "O" = "❌"; "P" = "❌"; "Q" = "✅"; "R" = "❌"; "S" = "❌"; "T" = "❌"; "U" = 2; "V" = 17; "W" = 11.8 };
    @{ "A" = "EC4255"; "B" = "Damini Solunke"; "C" = "daminisolunke@gmail.com"; "D" = "❌"; "E" = "❌"; "F" = "❌"; "G" = "❌"; "H" = "❌"; "I" = "❌"; "J" = "❌"; "K" = "❌"; "L" = "❌"; "M" = "❌"; "N" = "❌"; "O" = "❌"; "P" = "❌"; "Q" = "❌"; "R" = "❌"; "S" = "❌"; "T" = "❌"; "U" = 0; "V" = 17; "W" = 0 };
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 1
    $rowData = $data[$i]
    foreach ($c in $cols) {
        $cell = $ws.Range($c + $rowNum)
        $val = $rowData[$c]
        if (($rowNum -gt 1) -and ($numericCols -contains $c)) {
            $cell.Value = [double]$val
        } else {
            $cell.Value = $val
        }
    }
}
